# Fixing smaller issues in news:
#  1) "publications" sheet: turn numeric PubMed IDs (and empty PubMed cells)
#     in column I into PubMed link text, matching the style already used in
#     the DOI column.
#  2) "preprints" sheet: the "Pubmed" column (H) is always empty for this
#     sheet - remove it entirely, shifting DOI (old I) left into H.
#  3) "submissions" sheet: remove the unused "Journal" column (C), shifting
#     Date/Qualiperf/Authors Qualiperf/Projects left.

$wb = $excel.ActiveWorkbook

# --- 1) publications: column I (Pubmed) -----------------------------------
$pubs = $wb.Worksheets.Item("publications")
for ($r = 2; $r -le 34; $r++) {
    $cell = $pubs.Cells.Item($r, 9)
    $pmid = $cell.Text
    $cell.Value = '<a href="https://pubmed.ncbi.nlm.nih.gov/' + $pmid + '/">' + $pmid + '</a>'
}

# --- 2) preprints: drop the empty "Pubmed" column H ------------------------
$preprints = $wb.Worksheets.Item("preprints")
$preprints.Columns.Item(8).Delete()

# --- 3) submissions: drop the unused "Journal" column C --------------------
$submissions = $wb.Worksheets.Item("submissions")
$submissions.Columns.Item(3).Delete()
